# Update crypto price/volume data per latest scrape (GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.968.38"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.41%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.640.31"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.19%  "
$ws.Range("E4").Value = "  -0.24%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5092"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.26%  "
$ws.Range("E7").Value = "  -0.22%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2563"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.24%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06369"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.09%  "
$ws.Range("E10").Value = "  +0.29%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07746"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.53%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.279"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.49%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.646.90"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.38%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5450"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.61%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0₅7741"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.78%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.19"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.44%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "25.975.89"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.26%  "
$ws.Range("E18").Value = "  -0.22%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "196.14"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.24%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.424"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.18%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.923"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.043"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.36%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.004"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.17%  "
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "141.19"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.93%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1197"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +5.34%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.843"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.22%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.59"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.41%  "
$ws.Range("E29").Value = "  -0.29%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.04873"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.34%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.251"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.13%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.171"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.16%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.525"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.42%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.365"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.05%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.8932"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.65%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.144.26"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.71%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.577"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.11%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5444"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.29%  "
$ws.Range("E39").Value = "  -0.22%  "
$ws.Range("E40").Value = "  -0.28%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.522"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.99%  "
$ws.Range("E42").Value = "  +3.90%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8107"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.34%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.09"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.18%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.435"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.08%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.777.73"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.25%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4526"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.09%  "
$ws.Range("E48").Value = "  -0.73%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.9975"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.16%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05055"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.09%  "
$ws.Range("E51").Value = "  -0.39%  "

Write-Output "Applied cryptos update"
